$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 2
    3  = 0
    4  = 0
    5  = 0
    6  = 0
    7  = 3
    8  = 0
    9  = 0
    10 = 2
    11 = 0
    12 = 1
    13 = 2
    14 = 0
    15 = 1
    16 = 0
    17 = 2
    18 = 0
    19 = 0
    20 = 1
    21 = 1
    22 = 1
    23 = 2
    24 = 0
    25 = 0
    26 = 2
    27 = 1
    28 = 3
    29 = 0
    30 = 1
    31 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
